# Updated cryptos list on Thu Sep 19 23:28:23 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.040.08"
$ws.Range("E2").Value = "  +3.21%  "
$ws.Range("D3").Value = "2.472.82"
$ws.Range("E3").Value = "  +5.15%  "
$ws.Range("D5").Value = "'565.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "'142.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.45%  "
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "2.471.21"
$ws.Range("E9").Value = "  +5.21%  "
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("D11").Value = "'5.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "'0.351"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("D14").Value = "'26.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.26%  "
$ws.Range("D15").Value = "2.915.00"
$ws.Range("E15").Value = "  +5.20%  "
$ws.Range("D16").Value = "62.908.94"
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("D17").Value = "'0.0000141"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.30%  "
$ws.Range("D18").Value = "2.473.84"
$ws.Range("E18").Value = "  +6.02%  "
$ws.Range("D19").Value = "'11.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.52%  "
$ws.Range("D20").Value = "'340.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.70%  "
$ws.Range("E21").Value = "  +2.97%  "
$ws.Range("D22").Value = "'6.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'65.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  +4.83%  "
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("E29").Value = "  +6.76%  "
$ws.Range("D30").Value = "'6.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.04%  "
$ws.Range("D31").Value = "'1.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.17%  "
$ws.Range("D32").Value = "0.0₃0798"
$ws.Range("E32").Value = "  +7.72%  "
$ws.Range("D33").Value = "'176.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("E34").Value = "  +9.80%  "
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("D36").Value = "'18.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.75%  "
$ws.Range("D37").Value = "'372.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.27%  "
$ws.Range("B38").Value = "USDe"
$ws.Range("C38").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D38").Value = "'0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'4.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.82%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  +9.00%  "
$ws.Range("D42").Value = "'40.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.26%  "
$ws.Range("D43").Value = "'149.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.55%  "
$ws.Range("E44").Value = "  +4.39%  "
$ws.Range("D45").Value = "'20.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.42%  "
$ws.Range("E46").Value = "  +4.51%  "
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("D48").Value = "'0.0515"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("E49").Value = "  +4.09%  "
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("E51").Value = "  +3.91%  "
